# Apply odds updates to Jogos_da_Semana_FlashScore_2024-11-17.xlsx
# This script updates specific numeric (odds) cells across rows 2, 4, 5, 6, 10, 11, 13
# of the active worksheet, per the upstream data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2.38
$ws.Range("H2").Value = 2.88
$ws.Range("I2").Value = 3.4
$ws.Range("J2").Value = 3.25
$ws.Range("M2").Value = 1.14
$ws.Range("N2").Value = 5.5
$ws.Range("S2").Value = 1.67
$ws.Range("T2").Value = 2.1
$ws.Range("U2").Value = 2.38
$ws.Range("V2").Value = 1.53
$ws.Range("Z2").Value = 23
$ws.Range("AA2").Value = 26
$ws.Range("AH2").Value = 7
$ws.Range("AN2").Value = 4.33
$ws.Range("AT2").Value = 2.1
$ws.Range("AU2").Value = 10

# Row 4
$ws.Range("G4").Value = 1.65
$ws.Range("H4").Value = 3.5
$ws.Range("I4").Value = 6
$ws.Range("J4").Value = 2.3
$ws.Range("L4").Value = 6
$ws.Range("N4").Value = 8
$ws.Range("U4").Value = 2.1
$ws.Range("V4").Value = 1.67
$ws.Range("X4").Value = 7
$ws.Range("Z4").Value = 12
$ws.Range("AC4").Value = 8
$ws.Range("AE4").Value = 19
$ws.Range("AF4").Value = 67
$ws.Range("AH4").Value = 12
$ws.Range("AJ4").Value = 19
$ws.Range("AO4").Value = 9
$ws.Range("AV4").Value = 67
$ws.Range("AW4").Value = 7
$ws.Range("AZ4").Value = 126
$ws.Range("BA4").Value = 151

# Row 5
$ws.Range("G5").Value = 1.83
$ws.Range("H5").Value = 3.25
$ws.Range("I5").Value = 4.2
$ws.Range("J5").Value = 2.35
$ws.Range("K5").Value = 2.1
$ws.Range("L5").Value = 4.6
$ws.Range("N5").Value = 7.7
$ws.Range("O5").Value = 1.33
$ws.Range("P5").Value = 2.8
$ws.Range("Q5").Value = 1.98
$ws.Range("R5").Value = 1.65
$ws.Range("S5").Value = 1.39
$ws.Range("T5").Value = 2.57
$ws.Range("U5").Value = 1.8
$ws.Range("V5").Value = 1.8
$ws.Range("W5").Value = 6.6
$ws.Range("X5").Value = 8.5
$ws.Range("Z5").Value = 15.5
$ws.Range("AA5").Value = 15
$ws.Range("AB5").Value = 28
$ws.Range("AC5").Value = 8.5
$ws.Range("AD5").Value = 6.3
$ws.Range("AE5").Value = 15.5
$ws.Range("AH5").Value = 10.5
$ws.Range("AI5").Value = 23
$ws.Range("AJ5").Value = 14
$ws.Range("AK5").Value = 70
$ws.Range("AL5").Value = 45
$ws.Range("AM5").Value = 50
$ws.Range("AN5").Value = 3.65
$ws.Range("AO5").Value = 9
$ws.Range("AP5").Value = 17
$ws.Range("AQ5").Value = 30
$ws.Range("AS5").Value = 200
$ws.Range("AT5").Value = 2.52
$ws.Range("AU5").Value = 7.1
$ws.Range("AW5").Value = 6
$ws.Range("AX5").Value = 25

# Row 6
$ws.Range("G6").Value = 2.05
$ws.Range("I6").Value = 3.4
$ws.Range("J6").Value = 2.63
$ws.Range("L6").Value = 3.75
$ws.Range("U6").Value = 1.57
$ws.Range("V6").Value = 2.25
$ws.Range("W6").Value = 9.5
$ws.Range("X6").Value = 11
$ws.Range("Y6").Value = 9
$ws.Range("Z6").Value = 19
$ws.Range("AE6").Value = 12
$ws.Range("AG6").Value = 126
$ws.Range("AI6").Value = 19
$ws.Range("AJ6").Value = 12
$ws.Range("AK6").Value = 34
$ws.Range("AL6").Value = 23
$ws.Range("AN6").Value = 4.33
$ws.Range("AO6").Value = 11
$ws.Range("AX6").Value = 17

# Row 10
$ws.Range("Q10").Value = 1.57
$ws.Range("R10").Value = 2.35

# Row 11
$ws.Range("H11").Value = 2.9
$ws.Range("J11").Value = 4
$ws.Range("L11").Value = 3.5
$ws.Range("M11").Value = 1.13
$ws.Range("N11").Value = 6
$ws.Range("O11").Value = 1.53
$ws.Range("P11").Value = 2.38
$ws.Range("Q11").Value = 2.7
$ws.Range("R11").Value = 1.44
$ws.Range("S11").Value = 1.62
$ws.Range("T11").Value = 2.2
$ws.Range("U11").Value = 2.2
$ws.Range("V11").Value = 1.62
$ws.Range("AC11").Value = 6
$ws.Range("AF11").Value = 81
$ws.Range("AS11").Value = 351
$ws.Range("AT11").Value = 2.2
$ws.Range("AU11").Value = 9.5
$ws.Range("AX11").Value = 17
$ws.Range("AY11").Value = 34
$ws.Range("BA11").Value = 101

# Row 13
$ws.Range("O13").Value = 1.57
$ws.Range("P13").Value = 2.25
$ws.Range("U13").Value = 2.25
$ws.Range("V13").Value = 1.57
$ws.Range("AE13").Value = 21
$ws.Range("AH13").Value = 6.5
$ws.Range("BA13").Value = 126
